$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F (dSF) pulled from repulled source data.
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -1
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 6
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = -1
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 4
$ws.Range("F20").Value = 3
$ws.Range("F22").Value = 5
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = -4
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = -2
$ws.Range("F28").Value = -1
$ws.Range("F29").Value = 4
$ws.Range("F30").Value = -1
